$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.427221935424996
$ws.Range("B3").Value = 1.284323943879997
$ws.Range("B4").Value = 1.196362627868893
$ws.Range("B5").Value = 1.1604647620793
$ws.Range("B6").Value = 1.154500834163116
$ws.Range("B7").Value = 1.195878707289921
$ws.Range("B8").Value = 1.377998254199554
$ws.Range("B9").Value = 1.733274424020976
$ws.Range("B10").Value = 1.99304925807229
$ws.Range("B11").Value = 2.110934805366981
$ws.Range("B12").Value = 2.155531290526028
$ws.Range("B13").Value = 2.145928651782526
$ws.Range("B14").Value = 2.114604687648409
$ws.Range("B15").Value = 2.095412011724193
$ws.Range("B16").Value = 1.985339130512102
$ws.Range("B17").Value = 1.917737324343705
$ws.Range("B18").Value = 1.878827680503207
$ws.Range("B19").Value = 1.865649016799978
$ws.Range("B20").Value = 1.924936454307726
$ws.Range("B21").Value = 2.123806518121626
$ws.Range("B22").Value = 2.253520564414771
$ws.Range("B23").Value = 2.184314389785925
$ws.Range("B24").Value = 1.921681865934204
$ws.Range("B25").Value = 1.637373651876715
$ws.Range("C2").Value = 0.262848996876869
$ws.Range("C3").Value = 0.2468767195153134
$ws.Range("C4").Value = 0.2370509567479644
$ws.Range("C5").Value = 0.2330425290804214
$ws.Range("C6").Value = 0.2323766805954506
$ws.Range("C7").Value = 0.2369969148255393
$ws.Range("C8").Value = 0.2573458497254819
$ws.Range("C9").Value = 0.2970878127940182
$ws.Range("C10").Value = 0.3261718707782393
$ws.Range("C11").Value = 0.339375075837097
$ws.Range("C12").Value = 0.3443705491201285
$ws.Range("C13").Value = 0.3432948814173358
$ws.Range("C14").Value = 0.3397861446019306
$ws.Range("C15").Value = 0.3376363718325024
$ws.Range("C16").Value = 0.3253084295436111
$ws.Range("C17").Value = 0.3177383842711095
$ws.Range("C18").Value = 0.3133817488221382
$ws.Range("C19").Value = 0.3119062410688969
$ws.Range("C20").Value = 0.3185444946768712
$ws.Range("C21").Value = 0.3408168651346557
$ws.Range("C22").Value = 0.3553479485161404
$ws.Range("C23").Value = 0.3475948632912775
$ws.Range("C24").Value = 0.3181800662000001
$ws.Range("C25").Value = 0.2863555384331846
$ws.Range("D2").Value = 0.2253372448138506
$ws.Range("D3").Value = 0.2243814088732634
$ws.Range("D4").Value = 0.223867990835636
$ws.Range("D5").Value = 0.2236773157834406
$ws.Range("D6").Value = 0.2236467766638839
$ws.Range("D7").Value = 0.2238653441324061
$ws.Range("D8").Value = 0.2249924668718961
$ws.Range("D9").Value = 0.2277827930866607
$ws.Range("D10").Value = 0.2301829741705745
$ws.Range("D11").Value = 0.2313501960975657
$ws.Range("D12").Value = 0.2318029631245366
$ws.Range("D13").Value = 0.2317049739999391
$ws.Range("D14").Value = 0.2313872301718902
$ws.Range("D15").Value = 0.2311940027328063
$ws.Range("D16").Value = 0.2301082038728879
$ws.Range("D17").Value = 0.2294613533463945
$ws.Range("D18").Value = 0.2290964002958162
$ws.Range("D19").Value = 0.2289740548036576
$ws.Range("D20").Value = 0.2295294775707504
$ws.Range("D21").Value = 0.231480267642155
$ws.Range("D22").Value = 0.2328179373431709
$ws.Range("D23").Value = 0.2320982832269749
$ws.Range("D24").Value = 0.2294986570183397
$ws.Range("D25").Value = 0.2269661988934502
$ws.Range("F2").Value = 1.055527679690613
$ws.Range("F3").Value = 1.064485976268294
$ws.Range("F4").Value = 1.070665808268203
$ws.Range("F5").Value = 1.073354941408162
$ws.Range("F6").Value = 1.073811784726324
$ws.Range("F7").Value = 1.070701383343042
$ws.Range("F8").Value = 1.058475466339985
$ws.Range("F9").Value = 1.039893468838763
$ws.Range("F10").Value = 1.029532745360328
$ws.Range("F11").Value = 1.025534940208921
$ws.Range("F12").Value = 1.02412400060804
$ws.Range("F13").Value = 1.024423292089324
$ws.Range("F14").Value = 1.025416797801682
$ws.Range("F15").Value = 1.026038757150403
$ws.Range("F16").Value = 1.029808412457136
$ws.Range("F17").Value = 1.032304245264697
$ws.Range("F18").Value = 1.033807108945254
$ws.Range("F19").Value = 1.034327513796597
$ws.Range("F20").Value = 1.032031590977851
$ws.Range("F21").Value = 1.025122186654279
$ws.Range("F22").Value = 1.021206553545689
$ws.Range("F23").Value = 1.023241469785724
$ws.Range("F24").Value = 1.032154646246525
$ws.Range("F25").Value = 1.044342640518643
$ws.Range("G2").Value = 0.5091178494961497
$ws.Range("G3").Value = 0.5146137275156946
$ws.Range("G4").Value = 0.5184288881013117
$ws.Range("G5").Value = 0.5200942355342164
$ws.Range("G6").Value = 0.5203774424562013
$ws.Range("G7").Value = 0.5184508998158393
$ws.Range("G8").Value = 0.5109212643922376
$ws.Range("G9").Value = 0.4996600917148655
$ws.Range("G10").Value = 0.4935340957283216
$ws.Range("G11").Value = 0.4912159554734998
$ws.Range("G12").Value = 0.4904057108221735
$ws.Range("G13").Value = 0.4905772024946913
$ws.Range("G14").Value = 0.4911479405219552
$ws.Range("G15").Value = 0.4915063414690053
$ws.Range("G16").Value = 0.4936950383192311
$ws.Range("G17").Value = 0.4951578969553339
$ws.Range("G18").Value = 0.4960433901501773
$ws.Range("G19").Value = 0.4963507702748871
$ws.Range("G20").Value = 0.4949976079735166
$ws.Range("G21").Value = 0.4909784651153117
$ws.Range("G22").Value = 0.4887457499119137
$ws.Range("G23").Value = 0.4899012730356276
$ws.Range("G24").Value = 0.4950699360974866
$ws.Range("G25").Value = 0.5023302008476094
$ws.Range("H2").Value = 0.6614487526469262
$ws.Range("H3").Value = 0.6683050636529373
$ws.Range("H4").Value = 0.6728632893999134
$ws.Range("H5").Value = 0.6748084258694576
$ws.Range("H6").Value = 0.6751367067324949
$ws.Range("H7").Value = 0.67288916739966
$ws.Range("H8").Value = 0.6637404989705331
$ws.Range("H9").Value = 0.648564180911734
$ws.Range("H10").Value = 0.6390987371702863
$ws.Range("H11").Value = 0.6351583527060569
$ws.Range("H12").Value = 0.6337187863515226
$ws.Range("H13").Value = 0.6340264843579746
$ws.Range("H14").Value = 0.6350388650197587
$ws.Range("H15").Value = 0.6356658242148114
$ws.Range("H16").Value = 0.6393636064441921
$ws.Range("H17").Value = 0.6417256978934063
$ws.Range("H18").Value = 0.6431187115503079
$ws.Range("H19").Value = 0.6435962700298674
$ws.Range("H20").Value = 0.6414706887343726
$ws.Range("H21").Value = 0.634740077385004
$ws.Range("H22").Value = 0.6306476684464428
$ws.Range("H23").Value = 0.6328038206920326
$ws.Range("H24").Value = 0.6415858693035688
$ws.Range("H25").Value = 0.6523739052205144
$ws.Range("I2").Value = 0.6225890379057581
$ws.Range("I3").Value = 0.6336833811556257
$ws.Range("I4").Value = 0.6409301130735177
$ws.Range("I5").Value = 0.6439924292284793
$ws.Range("I6").Value = 0.6445075175213386
$ws.Range("I7").Value = 0.6409709704827495
$ws.Range("I8").Value = 0.6263240682000397
$ws.Range("I9").Value = 0.6010557018564544
$ws.Range("I10").Value = 0.584602455300983
$ws.Range("I11").Value = 0.5775771124644606
$ws.Range("I12").Value = 0.5749829636059864
$ws.Range("I13").Value = 0.5755387135722856
$ws.Range("I14").Value = 0.5773623627751459
$ws.Range("I15").Value = 0.5784880255695288
$ws.Range("I16").Value = 0.5850708332053003
$ws.Range("I17").Value = 0.5892269198889384
$ws.Range("I18").Value = 0.5916606158415636
$ws.Range("I19").Value = 0.5924920429178293
$ws.Range("I20").Value = 0.5887800226993107
$ws.Range("I21").Value = 0.5768249150960898
$ws.Range("I22").Value = 0.5693974695559945
$ws.Range("I23").Value = 0.5733262782484587
$ws.Range("I24").Value = 0.5889819269183505
$ws.Range("I25").Value = 0.6075210846492354
$ws.Range("J2").Value = 0.1768943251942074
$ws.Range("J3").Value = 0.1791500328625872
$ws.Range("J4").Value = 0.1806153593326338
$ws.Range("J5").Value = 0.1812326970291034
$ws.Range("J6").Value = 0.181336426124215
$ws.Range("J7").Value = 0.1806236031412585
$ws.Range("J8").Value = 0.1776554304714439
$ws.Range("J9").Value = 0.1724715797029379
$ws.Range("J10").Value = 0.1690503501220899
$ws.Range("J11").Value = 0.167577875106562
$ws.Range("J12").Value = 0.1670323357967298
$ws.Range("J13").Value = 0.1671492913586405
$ws.Range("J14").Value = 0.1675327516805121
$ws.Range("J15").Value = 0.1677692020058847
$ws.Range("J16").Value = 0.1691482669617912
$ws.Range("J17").Value = 0.17001575428921
$ws.Range("J18").Value = 0.1705226035329872
$ws.Range("J19").Value = 0.1706955698098271
$ws.Range("J20").Value = 0.1699225918651823
$ws.Range("J21").Value = 0.1674197929493282
$ws.Range("J22").Value = 0.1658543333303761
$ws.Range("J23").Value = 0.1666834206731878
$ws.Range("J24").Value = 0.1699646853061632
$ws.Range("J25").Value = 0.1738058678406915
$ws.Range("L2").Value = 0.4068172373015244
$ws.Range("L3").Value = 0.3956009180164273
$ws.Range("L4").Value = 0.3888257664023627
$ws.Range("L5").Value = 0.3860931927477651
$ws.Range("L6").Value = 0.3856411713809678
$ws.Range("L7").Value = 0.3887887987957583
$ws.Range("L8").Value = 0.4029268184140591
$ws.Range("L9").Value = 0.4315278708312746
$ws.Range("L10").Value = 0.4530644106797723
$ws.Range("L11").Value = 0.4629734198082076
$ws.Range("L12").Value = 0.4667415668328943
$ws.Range("L13").Value = 0.4659293297047782
$ws.Range("L14").Value = 0.4632831120404006
$ws.Range("L15").Value = 0.4616642794433261
$ws.Range("L16").Value = 0.4524190624023561
$ws.Range("L17").Value = 0.4467758929702939
$ws.Range("L18").Value = 0.4435406428405599
$ws.Range("L19").Value = 0.442447063390432
$ws.Range("L20").Value = 0.4473755279276332
$ws.Range("L21").Value = 0.464059943140569
$ws.Range("L22").Value = 0.4750562590509446
$ws.Range("L23").Value = 0.4691789829947481
$ws.Range("L24").Value = 0.4471044046509718
$ws.Range("L25").Value = 0.4236978361714563
$ws.Range("O2").Value = 2.306763694216016
$ws.Range("O3").Value = 2.332494378508471
$ws.Range("O4").Value = 2.349947261803479
$ws.Range("O5").Value = 2.357475048796942
$ws.Range("O6").Value = 2.358750122678998
$ws.Range("O7").Value = 2.350047101910562
$ws.Range("O8").Value = 2.315292107099083
$ws.Range("O9").Value = 2.260278754688997
$ws.Range("O10").Value = 2.227894638160222
$ws.Range("O11").Value = 2.214911782538536
$ws.Range("O12").Value = 2.210247397056634
$ws.Range("O13").Value = 2.211240743562939
$ws.Range("O14").Value = 2.214522988869561
$ws.Range("O15").Value = 2.216566282480315
$ws.Range("O16").Value = 2.22877832880522
$ws.Range("O17").Value = 2.236718271937164
$ws.Range("O18").Value = 2.241449699313577
$ws.Range("O19").Value = 2.243079932402566
$ws.Range("O20").Value = 2.235856014466378
$ws.Range("O21").Value = 2.213552072642159
$ws.Range("O22").Value = 2.200443867104696
$ws.Range("O23").Value = 2.20730542564479
$ws.Range("O24").Value = 2.236245322005885
$ws.Range("O25").Value = 2.27375208391112
